# data : case 1
# Replace the 32 x 2 numeric block (A1:B32) with the new computed values,
# and shrink the two data columns by ~1 character (autofit-style) to match
# the narrower recalculated numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.073001804238757018, 0.072900356218823958),
    @(-0.010914093097850852, 0.010574305989230481),
    @(0.092357461862057733, -0.09262059896749264),
    @(-0.19137326512317188, 0.19044991264627953),
    @(-0.18444991295203028, 0.18259298583518202),
    @(-0.073077217917661752, 0.073007617382399914),
    @(-0.053007617760810533, 0.052873946476738709),
    @(-0.068739831621981828, 0.068354148454440988),
    @(-0.062354148780749519, 0.062025729205593905),
    @(-0.056025729537026336, 0.055978332465009828),
    @(-0.051478332790683368, 0.051398772216998623),
    @(-0.04539877255051028, 0.04515282915074037),
    @(-0.039152829489757401, 0.03908581408959666),
    @(-0.027085814456018653, 0.027053371949834037),
    @(-0.021053372291694572, 0.02102783023169863),
    @(-0.01502783057471957, 0.015004530267366256),
    @(-0.0090045306118824442, 0.0089999996422251982),
    @(-0.051986263440451097, 0.051960568598431678),
    @(-0.042960568914011343, 0.042786563341294048),
    @(-0.018013837703884761, 0.018004301153256463),
    @(-0.0090043014738192539, 0.008999999679126347),
    @(-0.093931768167101026, 0.093623835335064243),
    @(-0.084623835653657054, 0.084124737255623749),
    @(-0.042124737721093375, 0.04199999953206035),
    @(-0.10351576823481423, 0.10337327479222225),
    @(-0.097373275110815172, 0.097195444925901597),
    @(-0.091195445246213147, 0.090608416638529476),
    @(-0.075244925565251641, 0.074523854558623981),
    @(-0.062523854916980426, 0.062166170904159657),
    @(-0.042166171299206745, 0.0420196121259786),
    @(-0.027019612502760637, 0.027000844026050075),
    @(-0.0060008444291952401, 0.0059999996611930229),
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Column A: 16.42578125 -> 15.42578125, Column B: 15.7109375 -> 14.7109375
# (each one character unit narrower). ColumnWidth is expressed in the
# "standard font character" units COM uses for Range/Columns.ColumnWidth.
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666
$ws.Columns.Item(2).ColumnWidth = 13.877604166666666
